$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.878.53'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '''1.629.22'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '''214.32'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '''0.503'
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '''19.61'
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").Value = '''0.0789'
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = '''1.854.41'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.23'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '''1.620.59'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '''62.72'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '''25.870.32'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").Value = '''0.999'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '''193.09'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '''4.38'
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("D22").Value = '''9.94'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = '''6.26'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Value = '''0.998'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").Value = '''142.13'
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = '''15.45'
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '''0.0500'
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").Value = '''2.42'
$ws.Range("E35").Value = '  +1.26%  '
$ws.Range("D36").Value = '''0.901'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").Value = '''1.131.27'
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''99.09'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''5.45'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").Value = '''1.764.09'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").Value = '''56.08'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("E48").Value = '  +3.73%  '
$ws.Range("D49").Value = '''1.46'
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").Value = '''7.61'
$ws.Range("E51").Value = '  +2.63%  '
